$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data rows that will be fully rewritten (old layout had rows 7-15; new layout rows 7-16)
$ws.Range("A7:D16").ClearContents()

# Safe (non-ambiguous) text/values - direct assignment
$ws.Range("A1").Value = 'Banque'
$ws.Range("A2").Value = 'Compte'
$ws.Range("A3").Value = 'Titulaire'
$ws.Range("B3").Value = 'ay 12s000 [or ]4sza.7e | U>Y994664173747/°2V/singhsonalt1@yBVMr SONAL'
$ws.Range("A4").Value = 'Période'
$ws.Range("B4").Value = '25/12/2019 - 23/12/2019'
$ws.Range("A7").Value = 'Date'
$ws.Range("B7").Value = 'Description'
$ws.Range("C7").Value = 'Montant'
$ws.Range("D7").Value = 'Sens'
$ws.Range("A8").Value = '25/12/2019'
$ws.Range("B8").Value = '‘Statement of Account: 2518XXXXXXXX0946 For Period: 26/9/2019 to'
$ws.Range("A9").Value = '23/12/2019'
$ws.Range("B9").Value = '|cr |88992 |218Z3SBIOC Ref No3000087061 ages 390.00 cw |ees.7= | UP7935119025907/"2V/sandeepchodhany@ybl/SANDEEP'
$ws.Range("D9").Value = 'Cr'
$ws.Range("B10").Value = '3 [DR [273.78 | ATMWDR 934616013196 +SINGHAL TOWER LABOUR'
$ws.Range("D10").Value = 'Dr'
$ws.Range("B11").Value = '[DR [3773.78 | IMPS-0UT/934616239787/SBIN00103160497052 ay 12s000 [or ]4sza.7e | U>Y994664173747/°2V/singhsonalt1@yBVMr SONAL Haazene 300000 [or |z3s:7= | UpY994857102800/P2v/e860521005@ybUMOHHAMMAD'
$ws.Range("D11").Value = 'Dr'
$ws.Range("B12").Value = '400 [cR [ _[2518001500010946:Int.Pd:01-09-2019 to 30-11-2019'
$ws.Range("D12").Value = 'Cr'
$ws.Range("B13").Value = '[DR [229.78 | UP/933884521318/P2M/IRCTCINAPP@ybI/IRCTC ouaENTD siooo Jor |7a7s7 | YPvaaae6sse04e0/"2V/7838767986@yb/Mr ABDUL'
$ws.Range("D13").Value = 'Dr'
$ws.Range("A14").Value = '21/11/2019'
$ws.Range("B14").Value = '[DR [227.97 | APBS REVRSAL DTD 01-06-2018'
$ws.Range("D14").Value = 'Dr'
$ws.Range("A15").Value = '15/11/2019'
$ws.Range("B15").Value = '200000 [DR [ | ATM WOR 931919032498 +SINGHAL TOWER LABOUR'
$ws.Range("D15").Value = 'Dr'
$ws.Range("A16").Value = '13/11/2019'
$ws.Range("B16").Value = '|cR [2,468.22 |21BZ3SBIOC Ref No3000072841'
$ws.Range("D16").Value = 'Cr'

# Ambiguous values that Excel would auto-convert to number/date - force text explicitly
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = '226.14'
$ws.Range("C9").ClearFormats()
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '12/12/2019'
$ws.Range("A10").ClearFormats()
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = '500.00'
$ws.Range("C10").ClearFormats()
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '12/12/2019'
$ws.Range("A11").ClearFormats()
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = '750.00'
$ws.Range("C11").ClearFormats()
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '09/12/2019'
$ws.Range("A12").ClearFormats()
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = '233.78'
$ws.Range("C12").ClearFormats()
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '04/12/2019'
$ws.Range("A13").ClearFormats()
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '508.19'
$ws.Range("C13").ClearFormats()
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '240.25'
$ws.Range("C14").ClearFormats()
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '468.22'
$ws.Range("C15").ClearFormats()
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '211.57'
$ws.Range("C16").ClearFormats()
